# Weekly update: insert the newest week's record for
# "Hortaliza, Vega Central Mapocho de Santiago - Perejil" as the new row 200,
# shifting all existing rows 200-236 down by one (to 201-237).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at 200; this pushes old rows 200..236 to 201..237
# and grows the used range to row 237.
$ws.Rows(200).Insert()

# Populate the newly inserted row 200 with this week's data.
$ws.Range("A200").Value = 9
$ws.Range("B200").Value = 'Vega Central Mapocho de Santiago'
$ws.Range("C200").Value = 'Metropolitana'
$ws.Range("D200").Value = 44505
$ws.Range("E200").Value = 13
$ws.Range("F200").Value = 100112044
$ws.Range("G200").Value = 'Perejil'
$ws.Range("H200").Value = 'Sin especificar'
$ws.Range("I200").Value = 'Primera'
$ws.Range("J200").Value = 106
$ws.Range("K200").Value = 8000
$ws.Range("L200").Value = 10000
$ws.Range("M200").Value = 9000
$ws.Range("N200").Value = '$/docena de atados'
$ws.Range("O200").Value = 'Región Metropolitana'
$ws.Range("P200").Value = 3000
$ws.Range("Q200").Value = 3
$ws.Range("R200").Value = 'Hortaliza'
